# Update the "想去人数" (want-to-go count) figures for the 2024-03-30 / 南宁
# update, matching output generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 9937
$wsExpo.Range("F3").Value = 218
$wsExpo.Range("F4").Value = 42
$wsExpo.Range("F5").Value = 591
$wsExpo.Range("F6").Value = 475

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 4

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 9937
$wsAll.Range("F3").Value = 218
$wsAll.Range("F4").Value = 42
$wsAll.Range("F5").Value = 591
$wsAll.Range("F6").Value = 4
$wsAll.Range("F7").Value = 475
